$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.718.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.601.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.28'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("E6").Value = '  -0.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.36%  '
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.826.26'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.603.43'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("E15").Value = '  +0.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.24'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.690.18'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("E18").Value = '  +1.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '211.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.63%  '
$ws.Range("E21").Value = '  +0.15%  '
$ws.Range("E22").Value = '  +0.92%  '
$ws.Range("E23").Value = '  +0.63%  '
$ws.Range("E24").Value = '  +1.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.59'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("E26").Value = '  +0.16%  '
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("E30").Value = '  +1.63%  '
$ws.Range("E32").Value = '  +1.60%  '
$ws.Range("E33").Value = '  +1.78%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.301.17'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.21%  '
$ws.Range("E35").Value = '  +0.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.608'
$ws.Range("D36").Style = "Normal"
$ws.Range("E37").Value = '  +1.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.16'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +20.52%  '
$ws.Range("E39").Value = '  -0.43%  '
$ws.Range("E40").Value = '  -1.68%  '
$ws.Range("E41").Value = '  -1.28%  '
$ws.Range("E42").Value = '  -0.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.781'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.21'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.737.13'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.03'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("E47").Value = '  -1.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0104'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.58%  '
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("E50").Value = '  +1.94%  '
$ws.Range("E51").Value = '  +0.00%  '
